$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sample number text "E7420" -> "E7420L" (shared by G2:G37)
$ws.Range("G2:G37").Value = "E7420L"

# Replace the "=FALSE()" formulas in H2:H37 with a literal boolean FALSE value
$ws.Range("H2:H37").Value = $false
